$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Comentarios")
$ws.Rows(10).Insert()
